$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number need to be forced to Text
# format first, so Excel stores them as strings (matching the original inline
# strings) instead of silently converting them to numeric values.

$ws.Range("D2").Value = '63.377.28'
$ws.Range("E2").Value = '  +2.18%  '
$ws.Range("D3").Value = '2.468.20'
$ws.Range("E3").Value = '  +1.88%  '
$ws.Range("E4").Value = '  +0.13%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '574.64'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +2.06%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '148.04'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +3.21%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E7").Value = '  -0.01%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.539'
$c.Style = "Normal"
$ws.Range("E8").Value = '  +1.67%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.113'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +3.54%  '
$ws.Range("E10").Value = '  +0.60%  '
$ws.Range("B11").Value = 'Toncoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '5.30'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +2.22%  '
$ws.Range("B12").Value = 'Cardano'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.361'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +3.16%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '27.05'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +3.65%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '0.0000183'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +5.62%  '
$ws.Range("D15").Value = '2.919.80'
$ws.Range("E15").Value = '  +2.13%  '
$ws.Range("D16").Value = '63.243.65'
$ws.Range("E16").Value = '  +2.04%  '
$ws.Range("D17").Value = '2.476.61'
$ws.Range("E17").Value = '  +1.81%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '11.46'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +1.61%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '7.32'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +7.44%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '328.29'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +1.44%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '4.22'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +2.17%  '
$ws.Range("B22").Value = 'SuiNetwork'
$ws.Range("C22").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '1.97'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +13.64%  '
$ws.Range("B23").Value = 'Dai'
$ws.Range("C23").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '0.998'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -0.14%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '67.18'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +0.07%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '629.90'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +13.57%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '8.85'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +0.77%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '0.0000105'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +13.00%  '
$ws.Range("D28").Value = '2.592.36'
$ws.Range("E28").Value = '  +1.93%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -0.07%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '8.38'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +2.36%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '0.146'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -0.66%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '1.91'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +2.11%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '5.17'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +9.02%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '1.53'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +1.64%  '
$ws.Range("E36").Value = '  -0.13%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.385'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +1.72%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '5.50'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +0.58%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '18.90'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +1.54%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '1.83'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +1.46%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '146.27'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -4.44%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '2.67'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +19.13%  '
$ws.Range("E43").Value = '  +0.16%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '150.07'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +1.88%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '3.76'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +3.52%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.0547'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +3.39%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '20.98'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +5.73%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.609'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +2.10%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '0.0238'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +4.61%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.0925'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +0.74%  '
$ws.Range("E51").Value = '  +0.76%  '
